# The "Warmup Plan" sheet had two leading rows (a "Properties" label row and
# a "Value" label row) above the real table header ("Phase", "Run", "Gmail",
# ...). Those two rows are removed so the real header becomes row 1 and the
# whole table shifts up by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")

$ws.Rows("1:2").Delete()

# Re-select the new header row (row 1), matching the row-header selection
# left behind in the saved file.
$ws.Rows("1:1").Select()
